$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (column C) date values for rows 2 through 15
# from 2023-10-05 (45204) to 2023-10-08 (45207), preserving the existing
# date serial number / formatting.
$oldDate = [datetime]::FromOADate(45204)

for ($row = 2; $row -le 15; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value() -eq $oldDate) {
        $cell.Value = 45207
    }
}
